$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 374, shifting rows 374:403 down to 375:404
$ws.Rows.Item(374).Insert()

# Populate the newly inserted row 374 with the new data record
$ws.Range("A374").Value = 8
$ws.Range("B374").Value = "Terminal La Palmera de La Serena"
$ws.Range("C374").Value = "Coquimbo"
$ws.Range("D374").Value = (Get-Date -Year 2023 -Month 4 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E374").Value = 4
$ws.Range("F374").Value = 100112021
$ws.Range("G374").Value = "Ají"
$ws.Range("H374").Value = "Inferno"
$ws.Range("I374").Value = "Primera"
$ws.Range("J374").Value = 500
$ws.Range("K374").Value = 12500
$ws.Range("L374").Value = 13000
$ws.Range("M374").Value = 12750
$ws.Range("N374").Value = "$/caja 15 kilos"
$ws.Range("O374").Value = "Provincia de Limarí"
$ws.Range("P374").Value = 850
$ws.Range("Q374").Value = 15
$ws.Range("R374").Value = "Hortaliza"
